# "changes to data table" - meanRate.xlsx / Sheet1
#
# 1. On Sheet1: move the per-row "Average" column from P to O (right next to
#    the data), add an "Average" header over O5:O6, insert a new "Average"
#    summary row (values pasted, not formulas) above the explanatory note,
#    and extend the 3-decimal scientific number format to the new cells.
# 2. Re-point the active tab from RVA to RVC.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Move the row-average formulas from column P to column O ------------
$ws.Range("O7").Formula = "=AVERAGE(D7:N7)"
$ws.Range("O8").Formula = "=AVERAGE(D8:N8)"
$ws.Range("O9").Formula = "=AVERAGE(D9:N9)"
$ws.Range("P7:P9").Clear()

# --- "Average" column header spanning the two header rows ---------------
$ws.Range("O5:O6").Merge()
$ws.Range("O5").Value = "Average"

# Extend the title banner merge (row 4) to cover the new column
$ws.Range("C4:N4").UnMerge()
$ws.Range("C4:O4").Merge()

# --- Insert a new row above the explanatory note for the column averages
$ws.Rows.Item(10).Insert()

$ws.Range("C10").Value = "Average"
$ws.Range("D10").Value = 0.017880266666666669
$ws.Range("E10").Value = 0.00102327
$ws.Range("F10").Value = 0.0021484333333333331
$ws.Range("G10").Value = 0.020665500000000003
$ws.Range("H10").Value = 0.0081763333333333341
$ws.Range("I10").Value = 0.0015666833333333333
$ws.Range("J10").Value = 0.0043837666666666662
$ws.Range("K10").Value = 0.010740333333333333
$ws.Range("L10").Value = 0.0020043999999999999
$ws.Range("M10").Value = 0.0019359666666666669
$ws.Range("N10").Value = 0.0054464000000000005
$ws.Range("O10").Formula = "=AVERAGE(O7:O9,D10:N10)"

# Extend the explanatory-note merge (now row 11) to cover the new column
$ws.Range("C11:N11").UnMerge()
$ws.Range("C11:O11").Merge()

# --- Apply the 3-decimal scientific number format to the data block -----
$ws.Range("D7:O10").NumberFormat = "0.000E+00"

# --- View bookkeeping -----------------------------------------------------
$ws.Range("Q11").Select()

$rva = $wb.Worksheets.Item("RVA")
$rvc = $wb.Worksheets.Item("RVC")
$rvc.Range("O14").Select()
$rvc.Activate()
